# Add a new "MainHud" row (row 10) to the Resource_Widget table and
# update the active selection to G10, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: A=10, B="MainHud" (new shared string), C=100, D=100,
# E=FALSE (boolean), F=0
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "MainHud"
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = 0

# Update dimension/selection to reflect the newly active cell
$ws.Range("G10").Select()
